$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "CustomLinkedList"
$ws.Range("C14").Value = "Impl of SLL with insert, delete, search, display etc"
$ws.Range("H14").Value = "CustomLinkedList"

$ws.Range("A15").Select()
